$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025/11/11"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "4.73"
$ws.Range("B2").Style = "Normal"

$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "2025/11/11"
$ws.Range("A8").Style = "Normal"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "7.48"
$ws.Range("B8").Style = "Normal"

$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "2025/11/11"
$ws.Range("A14").Style = "Normal"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "2.80"
$ws.Range("B14").Style = "Normal"

$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = "2025/11/11"
$ws.Range("A20").Style = "Normal"
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "12.07"
$ws.Range("B20").Style = "Normal"

$ws.Range("A26").NumberFormat = "@"
$ws.Range("A26").Value = "2025/11/11"
$ws.Range("A26").Style = "Normal"
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "9.93"
$ws.Range("B26").Style = "Normal"

$ws.Range("A32").NumberFormat = "@"
$ws.Range("A32").Value = "2025/11/11"
$ws.Range("A32").Style = "Normal"
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "24.76"
$ws.Range("B32").Style = "Normal"

$ws.Range("A38").NumberFormat = "@"
$ws.Range("A38").Value = "2025/11/11"
$ws.Range("A38").Style = "Normal"

$ws.Range("A44").NumberFormat = "@"
$ws.Range("A44").Value = "2025/11/11"
$ws.Range("A44").Style = "Normal"
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "11.53"
$ws.Range("B44").Style = "Normal"

$ws.Range("A50").NumberFormat = "@"
$ws.Range("A50").Value = "2025/11/11"
$ws.Range("A50").Style = "Normal"
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "11.64"
$ws.Range("B50").Style = "Normal"

$ws.Range("A56").NumberFormat = "@"
$ws.Range("A56").Value = "2025/11/11"
$ws.Range("A56").Style = "Normal"
$ws.Range("B56").NumberFormat = "@"
$ws.Range("B56").Value = "35.99"
$ws.Range("B56").Style = "Normal"

$ws.Range("A62").NumberFormat = "@"
$ws.Range("A62").Value = "2025/11/11"
$ws.Range("A62").Style = "Normal"
$ws.Range("B62").NumberFormat = "@"
$ws.Range("B62").Value = "11.89"
$ws.Range("B62").Style = "Normal"

$ws.Range("A68").NumberFormat = "@"
$ws.Range("A68").Value = "2025/11/11"
$ws.Range("A68").Style = "Normal"
$ws.Range("B68").NumberFormat = "@"
$ws.Range("B68").Value = "13.23"
$ws.Range("B68").Style = "Normal"

$ws.Range("A74").NumberFormat = "@"
$ws.Range("A74").Value = "2025/11/11"
$ws.Range("A74").Style = "Normal"
$ws.Range("B74").NumberFormat = "@"
$ws.Range("B74").Value = "16.25"
$ws.Range("B74").Style = "Normal"
